$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.470.86"
$ws.Range("E2").Value = "  -3.54%  "

$ws.Range("D3").Value = "3.397.20"
$ws.Range("E3").Value = "  -4.30%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.44%  "

$ws.Range("D8").Value = "3.394.02"
$ws.Range("E8").Value = "  -4.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.479"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.36%  "

$ws.Range("E10").Value = "  -10.88%  "

$ws.Range("E11").Value = "  -11.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.369"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -11.03%  "

$ws.Range("D13").Value = "3.972.38"
$ws.Range("E13").Value = "  -4.34%  "

$ws.Range("D14").Value = "3.447.02"
$ws.Range("E14").Value = "  -2.98%  "

$ws.Range("E15").Value = "  -11.48%  "

$ws.Range("E16").Value = "  -1.93%  "

$ws.Range("D17").Value = "64.497.93"
$ws.Range("E17").Value = "  -3.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -11.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -16.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.535"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.14%  "

$ws.Range("D26").Value = "3.530.63"

$ws.Range("E27").Value = "  -12.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  -12.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -13.83%  "

$ws.Range("D32").Value = "3.415.86"

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.09%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.141"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "170.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.41%  "

$ws.Range("E37").Value = "  -14.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -15.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.82%  "

$ws.Range("E40").Value = "  -14.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0750"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.791"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.63%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -16.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -12.36%  "

$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.23%  "

$ws.Range("E49").Value = "  -10.41%  "

$ws.Range("D50").Value = "2.177.35"
$ws.Range("E50").Value = "  -6.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.19%  "
